# Update section numbers in chapter title slides from Chapter 2/3 to Chapter 10.
$p = $ppt.ActivePresentation

# Slide 3: Title "3.1 Code" -> "10.1 Code"
$p.Slides.Item(3).Shapes.Item(1).TextFrame.TextRange.Text = "10.1 Code"

# Slide 4: Title "2.1 Code" -> "10.1 Code"
$p.Slides.Item(4).Shapes.Item(1).TextFrame.TextRange.Text = "10.1 Code"

# Slide 5: Title "3.2 Verify" -> "10.2 Verify"
$p.Slides.Item(5).Shapes.Item(1).TextFrame.TextRange.Text = "10.2 Verify"

# Slide 6: Title "3.2 Verify" -> "10.2 Verify"
$p.Slides.Item(6).Shapes.Item(1).TextFrame.TextRange.Text = "10.2 Verify"
